$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.825.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.22%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.628.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.01%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.67%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'214.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.504"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.57%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.14%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.0631"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.53%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'1.645.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.68%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'1.857.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.14%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'Polkadot"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'4.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.38%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0₃0753"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.13%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'Litecoin"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'62.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'25.849.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.53%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'192.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.79%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'9.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'6.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.83%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.82%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Monero"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'143.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.79%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'BinanceUSD"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.39%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.126"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.29%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E30").Value = "'  +0.05%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.30%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.02%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.94%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +2.09%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.898"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.74%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.136.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.74%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.547"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.60%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.02%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.58%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'99.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.58%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.800"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.52%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.80%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.767.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'Aave"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'56.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.41%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Cronos"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0527"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.01%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'RenderToken"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.20%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Mantle"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.415"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.34%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'7.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.34%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Algorand"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0960"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.48%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Updated cryptos list"
